$p = $ppt.ActivePresentation

# --- Slide 1 notes ---
$s = $p.Slides.Item(1)
$shp = $s.NotesPage.Shapes.Item(2)
$shp.TextFrame.TextRange.Text = "Good evening, we are the Spider Cartographers, and we’re going to explain how we chose to analyse and visualize transportation in the UK….so buckle up " + [char]10 + "" + [char]10 + "We will first begin with a presentation of the data sid, which will be followed up by an exploration of the website" + [char]10 + ""

# --- Slide 2 notes ---
$s = $p.Slides.Item(2)
$shp = $s.NotesPage.Shapes.Item(2)
$shp.TextFrame.TextRange.Text = "Ok what was our purpose? We want to gain a better understanding of how transport is experienced throughout the UK. Such an understanding can be used to inform both policy and the general public. " + [char]10 + "" + [char]10 + "As the UK shifts towards more sustainable modes of transport, our work aims to improve understanding of the current transportation ecosystem so as to better inform such a shift. " + [char]10 + "" + [char]10 + "Therefore, the aims of our analysis and website was to be able to identify and group different transport profiles across England and Wales at the MSOA level. These groups could then be mapped, along with transport flows, to see the geographic distribution of transport usage. " + [char]10 + "" + [char]10 + "We also demonstrate how we can use the identified profiles to see what demographic factors relate to the way in which we may use transport"

# --- Slide 3 notes ---
$s = $p.Slides.Item(3)
$shp = $s.NotesPage.Shapes.Item(2)
$shp.TextFrame.TextRange.Text = "Before identifying these transport profiles, there were a few steps we had to go through: gathering the data, cleaning it, transforming it and standardizing it"

# --- Slide 4 notes ---
$s = $p.Slides.Item(4)
$shp = $s.NotesPage.Shapes.Item(2)
$shp.TextFrame.TextRange.Text = "For this to work we required data on transportation usage across England and Wales. This was gathered at the MSOA  level given that this was the lowest geographical scale for which all data was available. The datasets covered transport access nodes, car ownership, commuter flow data between MSOAs and travel time, the latter being a proxy for accessibility. "

# --- Slide 5 notes ---
$s = $p.Slides.Item(5)
$shp = $s.NotesPage.Shapes.Item(2)
$shp.TextFrame.TextRange.Text = "The data that was received was in different formats and not all fit for our purpose. Some of the cleaning process we did included" + [char]10 + "Grouping the flows by origin MSOA and turning them into percentages so that they could be readily compared." + [char]10 + "Conducting a point in polygon analysis was for transport access nodes to see how many transport access nodes were in each MSOA. " + [char]10 + "Weighing travel times between MSOAs by the flows that occur between them. This allowed us to get an average of actual commuting times by mode for each MSOA, which was used as an accessibility metric"

# --- Slide 6 notes ---
$s = $p.Slides.Item(6)
$shp = $s.NotesPage.Shapes.Item(2)
$shp.TextFrame.TextRange.Text = "Once the data was cleaned, we had to transform and standardise it. Firstly transformation was used because extremes and outliers are likely influence cluster formation, especially for algorithms using distance-based metrics. Our variables were not skewed in the same degree or necessarily in the same direction, so we used two different transformation techniques" + [char]10 + "" + [char]10 + "Clustering results will also be affected by differences in units, ranges and variations. Therefore the data was standardised after transformation. Again due to the different variable distributions, no single standardization technique fit our data best, so we used three in order to compare the result." + [char]10 + "" + [char]10 + "The resulting outputs were all pushed to a single csv."
